# Renumber the TOC bookmarks / hyperlink anchors / PAGEREF fields after
# regenerating the table of contents.
#
# Mapping (old -> new):
#   _Toc18057388 -> _Toc40349825
#   _Toc18057389 -> _Toc40349826
#   _Toc18057390 -> _Toc40349827
#   _Toc18057391 -> _Toc40349828
#   _Toc18057392 -> _Toc40349829
#   _Toc18057393 -> _Toc40349830
#   _Toc18057394 -> _Toc40349831
#   _Toc18057395 -> _Toc40349832
#   _Toc18057396 -> _Toc40349833

$d = $word.ActiveDocument

$map = @{
    "_Toc18057388" = "_Toc40349825";
    "_Toc18057389" = "_Toc40349826";
    "_Toc18057390" = "_Toc40349827";
    "_Toc18057391" = "_Toc40349828";
    "_Toc18057392" = "_Toc40349829";
    "_Toc18057393" = "_Toc40349830";
    "_Toc18057394" = "_Toc40349831";
    "_Toc18057395" = "_Toc40349832";
    "_Toc18057396" = "_Toc40349833";
}

foreach ($old in $map.Keys) {
    $new = $map[$old]

    $rng = $d.Content
    $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null

    while ($rng.Find.Found) {
        $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
    }
}
